$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = "eman.tantawi@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"
$ws.Range("G4").Value = "Shimaa.ashraf@med.asu.edu.eg, aml.awwad@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg, Dr.hend.fariid@med.asu.edu.eg, aya.saeed@med.asu.edu.eg"
$ws.Range("G5").Value = "Sara_nabil@med.asu.edu.eg, wessam.atef@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg, aml.awwad@med.asu.edu.eg"
$ws.Range("G6").Value = "yassmina.fattoh@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, merna.said@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, maimustafa@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg"
$ws.Range("G7").Value = "mohamed.saleem@med.asu.edu.eg, afnan.fares@med.asu.edu.eg"
$ws.Range("G9").Value = "Mayarembaby@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg"
$ws.Range("G10").Value = "Rania.a.youssef@med.asu.edu.eg, System, Alshimaa_khaled@med.asu.edu.eg"
$ws.Range("G12").Value = "mohamed.saleem@med.asu.edu.eg, System"
$ws.Range("G14").Value = "nourhanmohamed@med.asu.edu.eg, marian.samir@med.asu.edu.eg"
$ws.Range("G15").Value = "nourhanmohamed@med.asu.edu.eg, marian.samir@med.asu.edu.eg, dr_ahmed_mostafa@med.asu.edu.eg, nourhan.mostafa@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg"
$ws.Range("G16").Value = "randa.rabea@med.asu.edu.eg, marian.samir@med.asu.edu.eg, Walaa.h.ghanima@med.asu.edu.eg, ManaratAleslam-Zefan@med.asu.edu.eg"
$ws.Range("G17").Value = "ola.m.abdelfattah@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, eman.samir@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg"
$ws.Range("G18").Value = "ola.m.abdelfattah@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, marina_atef@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg"
$ws.Range("G19").Value = "eman.tantawi@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"
$ws.Range("G21").Value = "Shimaa.ashraf@med.asu.edu.eg, aml.awwad@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg, Dr.hend.fariid@med.asu.edu.eg, aya.saeed@med.asu.edu.eg"
$ws.Range("G22").Value = "Sara_nabil@med.asu.edu.eg, wessam.atef@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg, aml.awwad@med.asu.edu.eg"
$ws.Range("G23").Value = "yassmina.fattoh@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, merna.said@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, maimustafa@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg"
$ws.Range("G24").Value = "mohamed.saleem@med.asu.edu.eg, afnan.fares@med.asu.edu.eg"
$ws.Range("G26").Value = "Mayarembaby@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg"
$ws.Range("G27").Value = "Rania.a.youssef@med.asu.edu.eg, System, Alshimaa_khaled@med.asu.edu.eg"
$ws.Range("G31").Value = "nourhanmohamed@med.asu.edu.eg, marian.samir@med.asu.edu.eg"
$ws.Range("G32").Value = "nourhanmohamed@med.asu.edu.eg, marian.samir@med.asu.edu.eg, dr_ahmed_mostafa@med.asu.edu.eg, nourhan.mostafa@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg"
$ws.Range("G33").Value = "randa.rabea@med.asu.edu.eg, marian.samir@med.asu.edu.eg, Walaa.h.ghanima@med.asu.edu.eg, ManaratAleslam-Zefan@med.asu.edu.eg"
$ws.Range("G34").Value = "ola.m.abdelfattah@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, eman.samir@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg"
$ws.Range("G35").Value = "Aya_hamed@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg"
$ws.Range("G36").Value = "eman.tantawi@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, System"
$ws.Range("G37").Value = "Kerelos.zareef@med.asu.edu.eg, NadaMohamed@med.asu.edu.eg, System"
$ws.Range("G40").Value = "yassmina.fattoh@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, merna.said@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, maimustafa@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg"
$ws.Range("G43").Value = "Mayarembaby@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg"
$ws.Range("G44").Value = "Rania.a.youssef@med.asu.edu.eg, System, Alshimaa_khaled@med.asu.edu.eg"
$ws.Range("G45").Value = "mohamed.saleem@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg, backup@backdoor.com, System"
$ws.Range("G46").Value = "mohamed.saleem@med.asu.edu.eg, System"
$ws.Range("G48").Value = "Aya.alaa-eldein@med.asu.edu.eg, marian.samir@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg"
$ws.Range("G50").Value = "nancy.abdelshafy@med.asu.edu.eg, Aya.alaa-eldein@med.asu.edu.eg, ManaratAleslam-Zefan@med.asu.edu.eg"
$ws.Range("G51").Value = "Salma.hassan@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, eman.samir@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg"
$ws.Range("G52").Value = "yasmin.m.senosy@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg"
$ws.Range("G53").Value = "eman.tantawi@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, System"
$ws.Range("G54").Value = "Kerelos.zareef@med.asu.edu.eg, NadaMohamed@med.asu.edu.eg, System"
$ws.Range("G57").Value = "yassmina.fattoh@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, merna.said@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, maimustafa@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg"
$ws.Range("G60").Value = "Mayarembaby@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg"
$ws.Range("G61").Value = "Rania.a.youssef@med.asu.edu.eg, System, Alshimaa_khaled@med.asu.edu.eg"
$ws.Range("G62").Value = "mohamed.saleem@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg, backup@backdoor.com, System"
$ws.Range("G65").Value = "Aya.alaa-eldein@med.asu.edu.eg, marian.samir@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg"
$ws.Range("G67").Value = "nancy.abdelshafy@med.asu.edu.eg, Aya.alaa-eldein@med.asu.edu.eg, ManaratAleslam-Zefan@med.asu.edu.eg"
$ws.Range("G68").Value = "Salma.hassan@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, eman.samir@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg"
$ws.Range("G69").Value = "yasmin.m.senosy@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg"
$ws.Range("G70").Value = "eman.tantawi@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"
$ws.Range("G72").Value = "Shimaa.ashraf@med.asu.edu.eg, wessam.atef@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg, Safa.hany@med.asu.edu.eg"
$ws.Range("G73").Value = "Shimaa.ashraf@med.asu.edu.eg, aml.awwad@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg, Dr.hend.fariid@med.asu.edu.eg, aya.saeed@med.asu.edu.eg"
$ws.Range("G74").Value = "Sara_nabil@med.asu.edu.eg, aml.awwad@med.asu.edu.eg, wessam.atef@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg, aya.saeed@med.asu.edu.eg"
$ws.Range("G75").Value = "maryam.ahmed@med.asu.edu.eg, Eman.m.elsaid@med.asu.edu.eg, merna.said@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg"
$ws.Range("G76").Value = "mohamed.saleem@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg"
$ws.Range("G77").Value = "Mayarembaby@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg"
$ws.Range("G78").Value = "Mayarembaby@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg, Alshimaa_khaled@med.asu.edu.eg"
$ws.Range("G79").Value = "Rania.a.youssef@med.asu.edu.eg, System, Alshimaa_khaled@med.asu.edu.eg"
$ws.Range("G83").Value = "Youstina.ibrahim@med.asu.edu.eg, marian.samir@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg"
$ws.Range("G84").Value = "Aya.alaa-eldein@med.asu.edu.eg, nancy.abdelshafy@med.asu.edu.eg, marian.samir@med.asu.edu.eg, ManaratAleslam-Zefan@med.asu.edu.eg"
$ws.Range("G85").Value = "ola.m.abdelfattah@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, nahla91@med.asu.edu.eg"
$ws.Range("G86").Value = "ola.m.abdelfattah@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, marina_atef@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg"
$ws.Range("G87").Value = "eman.tantawi@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"
$ws.Range("G89").Value = "Shimaa.ashraf@med.asu.edu.eg, aml.awwad@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg, Dr.hend.fariid@med.asu.edu.eg, aya.saeed@med.asu.edu.eg"
$ws.Range("G90").Value = "Sara_nabil@med.asu.edu.eg, aml.awwad@med.asu.edu.eg, wessam.atef@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg, aya.saeed@med.asu.edu.eg"
$ws.Range("G91").Value = "maryam.ahmed@med.asu.edu.eg, Eman.m.elsaid@med.asu.edu.eg, merna.said@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg"
$ws.Range("G92").Value = "Mayarembaby@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg"
$ws.Range("G93").Value = "Mayarembaby@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg"
$ws.Range("G94").Value = "Rania.a.youssef@med.asu.edu.eg, System, Alshimaa_khaled@med.asu.edu.eg"
$ws.Range("G98").Value = "nourhanmohamed@med.asu.edu.eg, user@user.com, Walaa.h.ghanima@med.asu.edu.eg, nourhanhosni@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg"
$ws.Range("G100").Value = "nourhanmohamed@med.asu.edu.eg, Aya.alaa-eldein@med.asu.edu.eg"
$ws.Range("G101").Value = "Aya.alaa-eldein@med.asu.edu.eg, nancy.abdelshafy@med.asu.edu.eg, marian.samir@med.asu.edu.eg, ManaratAleslam-Zefan@med.asu.edu.eg"
$ws.Range("G102").Value = "ola.m.abdelfattah@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, nahla91@med.asu.edu.eg"
$ws.Range("G103").Value = "ola.m.abdelfattah@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, gigi@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg"
$ws.Range("G104").Value = "eman.tantawi@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, System"
$ws.Range("G105").Value = "backup@backdoor.com, System"
$ws.Range("G106").Value = "Sara_nabil@med.asu.edu.eg, wessam.atef@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg, aml.awwad@med.asu.edu.eg"
$ws.Range("G107").Value = "Sara_nabil@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg, aml.awwad@med.asu.edu.eg"
$ws.Range("G108").Value = "esraa.mostafa@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, merna.said@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg"
$ws.Range("G109").Value = "mohamed.saleem@med.asu.edu.eg, afnan.fares@med.asu.edu.eg"
$ws.Range("G110").Value = "mohamed.saleem@med.asu.edu.eg, afnan.fares@med.asu.edu.eg, Alshimaa_khaled@med.asu.edu.eg"
$ws.Range("G111").Value = "Mayarembaby@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg"
$ws.Range("G112").Value = "Rania.a.youssef@med.asu.edu.eg, System, Alshimaa_khaled@med.asu.edu.eg"
$ws.Range("G114").Value = "mohamed.saleem@med.asu.edu.eg, System"
$ws.Range("G116").Value = "enas.omran@med.asu.edu.eg, nourhan.mostafa@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg"
$ws.Range("G117").Value = "user@user.com, enas.omran@med.asu.edu.eg, Amr-Saeed@med.asu.edu.eg, takyas@med.asu.edu.eg"
$ws.Range("G119").Value = "ola.m.abdelfattah@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, marinasorial@med.asu.edu.eg, marina_atef@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, eman.samir@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, nahla91@med.asu.edu.eg"
$ws.Range("G120").Value = "ola.m.abdelfattah@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, gigi@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg"
$ws.Range("G121").Value = "eman.tantawi@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, System"
$ws.Range("G122").Value = "backup@backdoor.com, System"
$ws.Range("G123").Value = "Sara_nabil@med.asu.edu.eg, wessam.atef@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg, aml.awwad@med.asu.edu.eg"
$ws.Range("G124").Value = "Sara_nabil@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg, aml.awwad@med.asu.edu.eg"
$ws.Range("G126").Value = "mohamed.saleem@med.asu.edu.eg, afnan.fares@med.asu.edu.eg"
$ws.Range("G127").Value = "mohamed.saleem@med.asu.edu.eg, afnan.fares@med.asu.edu.eg, Alshimaa_khaled@med.asu.edu.eg"
$ws.Range("G128").Value = "Mayarembaby@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg, Alshimaa_khaled@med.asu.edu.eg"
$ws.Range("G129").Value = "Rania.a.youssef@med.asu.edu.eg, System, Alshimaa_khaled@med.asu.edu.eg"
$ws.Range("G131").Value = "mohamed.saleem@med.asu.edu.eg, System"
$ws.Range("G133").Value = "enas.omran@med.asu.edu.eg, nourhan.mostafa@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg"
$ws.Range("G134").Value = "user@user.com, enas.omran@med.asu.edu.eg, Amr-Saeed@med.asu.edu.eg, takyas@med.asu.edu.eg"
$ws.Range("G136").Value = "ola.m.abdelfattah@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, marinasorial@med.asu.edu.eg, marina_atef@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, eman.samir@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, nahla91@med.asu.edu.eg"
$ws.Range("G137").Value = "ola.m.abdelfattah@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, gigi@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg"
$ws.Range("G138").Value = "eman.tantawi@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, System"
$ws.Range("G142").Value = "esraa.mostafa@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, merna.said@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg"
$ws.Range("G143").Value = "mohamed.saleem@med.asu.edu.eg, afnan.fares@med.asu.edu.eg"
$ws.Range("G144").Value = "Mayarembaby@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg, nourhanmamdouh@med.asu.edu.eg"
$ws.Range("G145").Value = "Mayarembaby@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg, Alshimaa_khaled@med.asu.edu.eg"
$ws.Range("G146").Value = "Mayarembaby@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg"
$ws.Range("G150").Value = "Youstina.ibrahim@med.asu.edu.eg, marian.samir@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg"
$ws.Range("G151").Value = "nourhanmohamed@med.asu.edu.eg, marian.samir@med.asu.edu.eg, hana.amr@med.asu.edu.eg, randa.rabea@med.asu.edu.eg, System"
